# "plot success rate 0-8": add a "success" column (D) next to list/count,
# flagging each row 0 (fail) / 1 (success) -- mirrors a pandas boolean
# column exported as text "0"/"1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell D1: same bold/boxed header style as B1/C1 ---
$ws.Range("D1").Value = "success"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data cells D2:D63: literal text "0"/"1" (not numbers) ---
$values = @("0","1","1","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0","0")

$rng = $ws.Range("D2:D63")
# Force text storage for these numeric-looking strings, then drop back to the
# default (unstyled) cell style so nothing extra sticks to the cells.
$rng.NumberFormat = "@"
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $values[$i]
}
$rng.Style = "Normal"
